$d = $word.ActiveDocument

# --- Merge the split runs back into single runs (retype-over via Find/Replace) ---
$d.Content.Find.Execute("MOUNTA.45 LASERCUT README", $true, $false, $false, $false, $false, $true, 1, $false, "MOUNTA.45 LASERCUT README", 2) | Out-Null
$d.Content.Find.Execute("Material: Delrin, Black", $true, $false, $false, $false, $false, $true, 1, $false, "Material: Delrin, Black", 2) | Out-Null
$d.Content.Find.Execute("Thickness: 0.193 inches", $true, $false, $false, $false, $false, $true, 1, $false, "Thickness: 0.193 inches", 2) | Out-Null

# --- Insert "Vendor: Ponoko" / blank / "File: ..." block between the two blank
#     lines that follow the title (i.e. before the 3rd paragraph) ---
$p3 = $d.Paragraphs.Item(3)

$p3.Range.InsertParagraphBefore()
$vendorPara = $d.Paragraphs.Item(3)
$vendorPara.Range.Text = "Vendor: Ponoko"
$vendorPara.Range.Font.NameFarEast = "Arial Unicode MS"
$vendorPara.Range.Font.NameBi = "Arial Unicode MS"

$p3 = $d.Paragraphs.Item(4)
$p3.Range.InsertParagraphBefore()

$p3 = $d.Paragraphs.Item(5)
$p3.Range.InsertParagraphBefore()
$filePara = $d.Paragraphs.Item(5)
$filePara.Range.Text = "File: MOUNTA.45_P2_0.193BLACK_DELRIN_40.eps"
$filePara.Range.Font.NameFarEast = "Arial Unicode MS"
$filePara.Range.Font.NameBi = "Arial Unicode MS"

# --- Append a blank line and a "Yield: 40" paragraph after "Thickness..." ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$yieldPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$yieldPara.Range.Text = "Yield: 40"
$yieldPara.Range.Font.NameFarEast = "Arial Unicode MS"
$yieldPara.Range.Font.NameBi = "Arial Unicode MS"
